$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the spicule type label in A1 from "tylostyle" to "Oxea"
$ws.Range("A1").Value = "Oxea"
